# DbSchema.xlsx update
# - rename the three table-header cells to their lowercase snake_case table names
# - remove the obsolete "session_date" row from the SessionDetails block
# - update the saved selection and page setup to match the resaved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename table name headers (order matters for shared-string table layout)
$ws.Range("A5").Value = "courses"
$ws.Range("A1").Value = "instructors"
$ws.Range("A16").Value = "session_details"

# Drop the "session_date" row from the SessionDetails table; rows below shift up
$ws.Rows(20).Delete()

# Match the active selection recorded in the saved file
$ws.Range("D16").Select()

# Page setup as captured on save
$ws.PageSetup.PaperSize = 512
$ws.PageSetup.Orientation = 1
